{"js": "// Office.js (Word JavaScript API) script implementing the target diff.\n// Body of: async (context) => { ... }\n\n// ---------------------------------------------------------------------\n// 1) Delete the paragraph \"At the end of this practical, upload this\n//    document to the Brightspace assignment\" entirely (it is removed in\n//    the target, and the following \"1. Characterize weather conditions\"\n//    heading takes its place).\n// ---------------------------------------------------------------------\n{\n  const results = context.document.body.search(\n    \"At the end of this practical, upload this document to the Brightspace assignment\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const para = results.items[0].paragraphs.getFirst();\n    para.delete();\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 2) \"... per day (or group of days). Insert graphs ...\" ->\n//    \"... per day (or group of days) (similar as you did in step 1).\n//    Insert graphs ...\"\n// ---------------------------------------------------------------------\n{\n  const results = context.document.body.search(\"or group of days)\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\" (similar as you did in step 1)\", \"After\");\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 3) \"... or make an\" -> \"... or make an additional/different table\"\n// ---------------------------------------------------------------------\n{\n  const results = context.document.body.search(\n    \"replace the \\u2018...\\u2019 in the table below the variable you choose). If needed extend the table, or make an\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\" additional/different table\", \"After\");\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 4) \"Describe your findings in concisely.\" -> \"Describe your findings\n//    concisely below.\"\n// ---------------------------------------------------------------------\n{\n  const results = context.document.body.search(\"findings in concisely.\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"findings concisely below.\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 5) \"Crop factors and reference evapotranspiration\" (heading 6) ->\n//    \"Crop factors and reference evapotranspiration for bare soil\"\n// ---------------------------------------------------------------------\n{\n  const results = context.document.body.search(\"Crop factors and reference evapotranspiration\", {\n    matchCase: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\" for bare soil\", \"After\");\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 6) \"What is the overall magnitude the crop factor? Is the crop factor\n//    constant...\" -> \"...the crop factor for the bare soil field? Is the\n//    crop factor constant...\"\n// ---------------------------------------------------------------------\n{\n  const results = context.document.body.search(\"What is the overall magnitude the crop factor?\", {\n    matchCase: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"What is the overall magnitude the crop factor for the bare soil field?\",\n      \"Replace\"\n    );\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 7) Remove the \"Site\" table header text (keep the empty, still-bold\n//    paragraph in the cell).\n// ---------------------------------------------------------------------\n{\n  const results = context.document.body.search(\"Site\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 8) \"Typical values\" -> \"Typical value of the crop factor (mean,\n//    median, \\u2026)\"\n// ---------------------------------------------------------------------\n{\n  const results = context.document.body.search(\"Typical values\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"Typical value of the crop factor (mean, median, \\u2026)\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 9) \"Variability\" -> \"Variability of crop factor from day to day (how\n//    much, how, when)\"\n// ---------------------------------------------------------------------\n{\n  const results = context.document.body.search(\"Variability\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"Variability of crop factor from day to day (how much, how, when)\",\n      \"Replace\"\n    );\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 10) \"What determines the variabililty\" -> \"What determines the\n//     day-to-day variation?\"\n// ---------------------------------------------------------------------\n{\n  const results = context.document.body.search(\"What determines the variabililty\", {\n    matchCase: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"What determines the day-to-day variation?\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) script implementing the target diff.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# Unicode helpers for the curly quotes used in the source text.\n$lsq = [char]0x2018   # '\n$rsq = [char]0x2019   # '\n$hellip = [char]0x2026 # \u2026\n\n# ---------------------------------------------------------------------\n# 1) Delete the paragraph \"At the end of this practical, upload this\n#    document to the Brightspace assignment\" entirely (it is removed in\n#    the target, and the following \"1. Characterize weather conditions\"\n#    heading takes its place).\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$found = $rng.Find.Execute(\"At the end of this practical, upload this document to the Brightspace assignment\")\nif ($found) {\n    $rng.Paragraphs(1).Range.Delete()\n}\n\n# ---------------------------------------------------------------------\n# 2) \"... per day (or group of days). Insert graphs ...\" ->\n#    \"... per day (or group of days) (similar as you did in step 1).\n#    Insert graphs ...\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\n    \"or group of days). \",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"or group of days) (similar as you did in step 1). \",\n    2\n) | Out-Null\n\n# ---------------------------------------------------------------------\n# 3) \"... or make an\" -> \"... or make an additional/different table\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$findText = \"replace the \" + $lsq + \"...\" + $rsq + \" in the table below the variable you choose). If needed extend the table, or make an\"\n$replaceText = $findText + \" additional/different table\"\n$rng.Find.Execute(\n    $findText,\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    $replaceText,\n    2\n) | Out-Null\n\n# ---------------------------------------------------------------------\n# 4) \"Describe your findings in concisely.\" -> \"Describe your findings\n#    concisely below.\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\n    \"findings in concisely.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"findings concisely below.\",\n    2\n) | Out-Null\n\n# ---------------------------------------------------------------------\n# 5) \"Crop factors and reference evapotranspiration\" (heading 6) ->\n#    \"Crop factors and reference evapotranspiration for bare soil\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\n    \"Crop factors and reference evapotranspiration\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Crop factors and reference evapotranspiration for bare soil\",\n    2\n) | Out-Null\n\n# ---------------------------------------------------------------------\n# 6) \"What is the overall magnitude the crop factor? Is the crop factor\n#    constant...\" -> \"...the crop factor for the bare soil field? Is the\n#    crop factor constant...\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\n    \"What is the overall magnitude the crop factor?\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"What is the overall magnitude the crop factor for the bare soil field?\",\n    2\n) | Out-Null\n\n# ---------------------------------------------------------------------\n# 7) Remove the \"Site\" table header text (keep the empty, still-bold\n#    paragraph in the cell).\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Site\")\nif ($found) {\n    $rng.Text = \"\"\n}\n\n# ---------------------------------------------------------------------\n# 8) \"Typical values\" -> \"Typical value of the crop factor (mean,\n#    median, \u2026)\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\n    \"Typical values\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Typical value of the crop factor (mean, median, \" + $hellip + \")\",\n    2\n) | Out-Null\n\n# ---------------------------------------------------------------------\n# 9) \"Variability\" -> \"Variability of crop factor from day to day (how\n#    much, how, when)\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\n    \"Variability\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Variability of crop factor from day to day (how much, how, when)\",\n    2\n) | Out-Null\n\n# ---------------------------------------------------------------------\n# 10) \"What determines the variabililty\" -> \"What determines the\n#     day-to-day variation?\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\n    \"What determines the variabililty\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"What determines the day-to-day variation?\",\n    2\n) | Out-Null\n"}
